# Update the "Estado de Cuenta" worksheet with the newly reconciled
# "Valor Mora" amounts for the MANUELA HINCAPIE BOTERO records.
#
# Periodo Mora 2005 (row 17) now carries the 33125 amount, while
# Periodo Mora 1910 (row 24) now carries the 20979 amount (the two
# values were swapped as part of the database refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F17").Value = 33125
$ws.Range("F24").Value = 20979
